# Update the Notes text for "5. Longest Palindromic Substring" (row 27, column D)
# with the expanded explanation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D27").Value = "For palindrome, we take each index as a midpoint and expand outwards. We use a dp[n][n] array to track palindromes starting at i,j, and also a variable for palindrome start and max length. At the end we return s.substring(palindrome_starts_at, palindrome_starts_at + max_len);"

# Add a new row to the table for "230. Kth Smallest Element in a BST"
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

$ws.Range("A29").Value = "230. Kth Smallest Element in a BST"

$ws.Range("B29").Value = "Medium"
$ws.Range("B29").Interior.Color = $ws.Range("B28").Interior.Color()

$ws.Range("C29").Value = "Trees"

$ws.Range("E29").Value = "https://leetcode.com/problems/kth-smallest-element-in-a-bst/solutions/63783/two-easiest-in-order-traverse-java/ "
$ws.Hyperlinks.Add($ws.Range("E29"), "https://leetcode.com/problems/kth-smallest-element-in-a-bst/solutions/63783/two-easiest-in-order-traverse-java/")
$ws.Range("E29").Style = $ws.Range("E28").Style()

$ws.Range("D29").Value = "As BST is already sorted, we can visit all nodes and store in an array to solve, optimally, the BST is already sorted so no array is needed. Use a stack and solve iteratively to immediately return the kth smallest node."

# Update the view: scroll so row 7 is near the top and select C37 (matches author's
# last on-screen state before saving).
$ws.Range("C37").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "done"
